# "Generate Report for Handback"
# Records the handback result for the 71fb0c0c-977c-4368-b211-1a5edd9535e1
# source file on both the zh-cn and de-de sheets: the returned xlf did not
# match the latest handoff, so we stamp a handback datetime + error detail,
# and point "Latest Target File" back at the source .md (with a hyperlink
# matching the style already used for column A).

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/844fc8cec6e56c7629c91935d38737c19c02ab34/e2e/71fb0c0c-977c-4368-b211-1a5edd9535e1.md"
$targetMdName = "71fb0c0c-977c-4368-b211-1a5edd9535e1.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/628d847701d8d7b43b8c10dad0f57a3c32968c47/e2e/71fb0c0c-977c-4368-b211-1a5edd9535e1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/844fc8cec6e56c7629c91935d38737c19c02ab34/e2e/71fb0c0c-977c-4368-b211-1a5edd9535e1.md."

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Latest Target File" / "Latest Handback File" / "Error Detail"
# columns to match the other wide (40) columns, same as column A.
$wideWidth = $wsZh.Columns.Item(1).ColumnWidth
$wsZh.Columns.Item(9).ColumnWidth = $wideWidth
$wsZh.Columns.Item(10).ColumnWidth = $wideWidth
$wsZh.Columns.Item(16).ColumnWidth = $wideWidth

# Row 5 == 71fb0c0c-977c-4368-b211-1a5edd9535e1
$wsZh.Range("J5").Value = $wsZh.Range("G5").Value2
$wsZh.Range("K5").Value = "2016-09-05 10:27:10"
$wsZh.Range("P5").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $targetMdUrl, [Type]::Missing, [Type]::Missing, $targetMdName) | Out-Null

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wideWidth = $wsDe.Columns.Item(1).ColumnWidth
$wsDe.Columns.Item(9).ColumnWidth = $wideWidth
$wsDe.Columns.Item(10).ColumnWidth = $wideWidth
$wsDe.Columns.Item(16).ColumnWidth = $wideWidth

# Row 5 == 71fb0c0c-977c-4368-b211-1a5edd9535e1
$wsDe.Range("J5").Value = $wsDe.Range("G5").Value2
$wsDe.Range("K5").Value = "2016-09-05 10:27:30"
$wsDe.Range("P5").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $targetMdUrl, [Type]::Missing, [Type]::Missing, $targetMdName) | Out-Null
